# Auto-generated PowerShell COM-interop script.
# Applies updated crypto price/volume figures (and the two re-ranked row swaps:
# Polkadot/Polygon at rows 14-15, and Algorand/BabyDogeCoin at rows 49-50)
# as captured by the authoritative XML diff for this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.352.62"
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = "'1.652.03"
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = "'213.36"
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = "'23.70"
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').Value = "'0.0877"
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = "'1.885.00"
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = "'1.647.82"
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.572"
$ws.Range('E14').Value = '  +4.14%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = "'4.06"
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = "'65.65"
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = "'27.344.65"
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = "'232.24"
$ws.Range('E18').Value = '  -6.45%  '
$ws.Range('D19').Value = "'0.0₃0725"
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = "'7.51"
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').Value = "'4.37"
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').Value = "'9.18"
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').Value = "'146.86"
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').Value = "'7.17"
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = "'15.88"
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').Value = "'1.458.06"
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('D35').Value = "'1.55"
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').Value = "'0.908"
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = "'0.572"
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = "'5.48"
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').Value = "'65.12"
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = "'1.792.83"
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('D46').Value = "'0.785"
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('D48').Value = "'88.10"
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.101"
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = "'0.0₇0999"
$ws.Range('E50').Value = '  -9.43%  '
$ws.Range('D51').Value = "'7.75"
$ws.Range('E51').Value = '  +0.33%  '
